$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff (cryptos list refresh).
# A leading apostrophe forces text interpretation so Excel doesn't
# auto-convert numeric-looking strings (e.g. "589.29") to numbers,
# preserving the original General/text cell formatting.
$ws.Range("D2").Value = "'69.371.19"
$ws.Range("E2").Value = "'  +2.41%  "
$ws.Range("D3").Value = "'3.387.76"
$ws.Range("E3").Value = "'  +1.66%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'589.29"
$ws.Range("E5").Value = "'  +1.45%  "
$ws.Range("D6").Value = "'180.89"
$ws.Range("E6").Value = "'  +3.05%  "
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E8").Value = "'  +0.92%  "
$ws.Range("D9").Value = "'0.197"
$ws.Range("E9").Value = "'  +8.82%  "
$ws.Range("E10").Value = "'  +1.53%  "
$ws.Range("D11").Value = "'48.72"
$ws.Range("E11").Value = "'  +4.79%  "
$ws.Range("D12").Value = "'0.0000288"
$ws.Range("E12").Value = "'  +6.16%  "
$ws.Range("D13").Value = "'687.65"
$ws.Range("E13").Value = "'  -2.41%  "
$ws.Range("E14").Value = "'  +2.24%  "
$ws.Range("D15").Value = "'3.940.24"
$ws.Range("E15").Value = "'  +1.60%  "
$ws.Range("D16").Value = "'69.449.96"
$ws.Range("E16").Value = "'  +2.54%  "
$ws.Range("B17").Value = "'TRON"
$ws.Range("C17").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.120"
$ws.Range("E17").Value = "'  +1.87%  "
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.341.99"
$ws.Range("E18").Value = "'  +0.29%  "
$ws.Range("E19").Value = "'  +2.24%  "
$ws.Range("D20").Value = "'11.43"
$ws.Range("E20").Value = "'  +3.76%  "
$ws.Range("E21").Value = "'  +1.02%  "
$ws.Range("E22").Value = "'  -0.12%  "
$ws.Range("D23").Value = "'17.12"
$ws.Range("E23").Value = "'  +0.85%  "
$ws.Range("D24").Value = "'104.51"
$ws.Range("E24").Value = "'  +6.07%  "
$ws.Range("E25").Value = "'  +1.73%  "
$ws.Range("E26").Value = "'  +1.51%  "
$ws.Range("D27").Value = "'9.64"
$ws.Range("E27").Value = "'  +1.58%  "
$ws.Range("D28").Value = "'34.48"
$ws.Range("E28").Value = "'  +4.05%  "
$ws.Range("D29").Value = "'8.71"
$ws.Range("E29").Value = "'  +1.83%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "'  -1.37%  "
$ws.Range("D31").Value = "'11.22"
$ws.Range("E31").Value = "'  +2.10%  "
$ws.Range("D32").Value = "'557.92"
$ws.Range("E32").Value = "'  -2.55%  "
$ws.Range("D33").Value = "'3.65"
$ws.Range("E33").Value = "'  +10.10%  "
$ws.Range("D34").Value = "'0.107"
$ws.Range("E34").Value = "'  +1.12%  "
$ws.Range("D35").Value = "'58.11"
$ws.Range("E35").Value = "'  +1.01%  "
$ws.Range("E36").Value = "'  +0.12%  "
$ws.Range("D37").Value = "'3.708.15"
$ws.Range("E37").Value = "'  -0.09%  "
$ws.Range("D38").Value = "'0.140"
$ws.Range("E38").Value = "'  +6.81%  "
$ws.Range("D39").Value = "'35.02"
$ws.Range("E39").Value = "'  +2.94%  "
$ws.Range("B40").Value = "'PEPE"
$ws.Range("C40").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "'0.0₃0710"
$ws.Range("E40").Value = "'  +5.11%  "
$ws.Range("B41").Value = "'Stacks"
$ws.Range("C41").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.25"
$ws.Range("E41").Value = "'  +1.63%  "
$ws.Range("D42").Value = "'2.70"
$ws.Range("E42").Value = "'  +2.30%  "
$ws.Range("D43").Value = "'0.340"
$ws.Range("E43").Value = "'  +1.43%  "
$ws.Range("D44").Value = "'0.0420"
$ws.Range("E44").Value = "'  +3.19%  "
$ws.Range("E45").Value = "'  +0.07%  "
$ws.Range("D46").Value = "'2.66"
$ws.Range("E46").Value = "'  -0.67%  "
$ws.Range("E47").Value = "'  +1.17%  "
$ws.Range("D48").Value = "'1.39"
$ws.Range("E48").Value = "'  +4.99%  "
$ws.Range("E49").Value = "'  -0.04%  "
$ws.Range("D50").Value = "'133.07"
$ws.Range("E50").Value = "'  +3.21%  "
$ws.Range("D51").Value = "'2.60"
$ws.Range("E51").Value = "'  -2.19%  "
